$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# The only real content change: cell E8's text changes from "Good Morning"
# to "GIT UPDATE" (all the shared-string index churn in the diff is just
# a consequence of "Good Morning" being dropped from the shared strings
# table and "GIT UPDATE" being appended at the end).
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the last active cell selection recorded in the saved file.
$ws.Range("E8").Select()
